$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 08:33"

# Swap country labels for rows 123/124 (Taiwan <-> Guinea Ecuatorial)
$ws.Range("A123").Value = "Taiwan"
$ws.Range("A124").Value = "Guinea Ecuatorial"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1263224
$ws.Range("C4").Value = 132
$ws.Range("E4").Value = 975306
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 74809

# Row 66 - Ghana
$ws.Range("D66").Value = 303
$ws.Range("E66").Value = 2770

# Row 114 - Georgia
$ws.Range("B114").Value = 615
$ws.Range("C114").Value = 5
$ws.Range("D114").Value = 275
$ws.Range("E114").Value = 331

# Row 123 - now Taiwan
$ws.Range("B123").Value = 440
$ws.Range("C123").Value = 1
$ws.Range("D123").Value = 347
$ws.Range("E123").Value = 87
$ws.Range("H123").Value = 6

# Row 124 - now Guinea Ecuatorial
$ws.Range("D124").Value = 13
$ws.Range("E124").Value = 422
$ws.Range("H124").Value = 4
